$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 594.75
$ws.Range("I17").Value = 200
$ws.Range("J17").Value = 615.5263
$ws.Range("K17").Value = 600
$ws.Range("L17").Value = 1846.5789
$ws.Range("M17").Value = -432
$ws.Range("N17").Value = -2182.5789

$ws.Range("H33").Value = 23331.727
$ws.Range("I33").Value = 25653.95
$ws.Range("J33").Value = 109.5
$ws.Range("K33").Value = 25653.95
$ws.Range("L33").Value = 109.5
$ws.Range("M33").Value = -25424.95
$ws.Range("N33").Value = -567.5

$ws.Range("H112").Value = 2517.5
$ws.Range("J112").Value = 2774.2856
$ws.Range("L112").Value = 8322.856800000001
$ws.Range("N112").Value = -10538.8568

$ws.Range("H132").Value = 44332.582
$ws.Range("I132").Value = 58238.832
$ws.Range("K132").Value = 174716.496
$ws.Range("M132").Value = -172186.496

$ws.Range("H137").Value = 12501313
$ws.Range("I137").Value = 22501036
$ws.Range("J137").Value = 1659.5625
$ws.Range("K137").Value = 67503108
$ws.Range("L137").Value = 4978.6875
$ws.Range("M137").Value = -67500558
$ws.Range("N137").Value = -10078.6875

$ws.Range("H138").Value = 11931.4
$ws.Range("I138").Value = 19700
$ws.Range("J138").Value = 9989.25
$ws.Range("K138").Value = 59100
$ws.Range("L138").Value = 29967.75
$ws.Range("M138").Value = -53960
$ws.Range("N138").Value = -40247.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5623.625
$ws.Range("J2").Value = 4362
$ws.Range("L2").Value = 4362
$ws.Range("N2").Value = -4588

$ws.Range("H32").Value = 3425.2666
$ws.Range("I32").Value = 3536.0715
$ws.Range("K32").Value = 3536.0715
$ws.Range("M32").Value = -3249.0715

$ws.Range("H74").Value = 650322.4399999999
$ws.Range("I74").Value = 4343.027
$ws.Range("J74").Value = 4633862
$ws.Range("K74").Value = 4343.027
$ws.Range("L74").Value = 4633862
$ws.Range("M74").Value = -3469.027
$ws.Range("N74").Value = -4635610

$ws.Range("H77").Value = 650322.4399999999
$ws.Range("I77").Value = 4343.027
$ws.Range("J77").Value = 4633862
$ws.Range("K77").Value = 21715.135
$ws.Range("L77").Value = 23169310
$ws.Range("M77").Value = -17347.135
$ws.Range("N77").Value = -23178046

$ws.Range("H116").Value = 5623.625
$ws.Range("J116").Value = 4362
$ws.Range("L116").Value = 4362
$ws.Range("N116").Value = -8950

$ws.Range("H122").Value = 6674.6665
$ws.Range("I122").Value = 5972.75
$ws.Range("K122").Value = 17918.25
$ws.Range("M122").Value = -15468.25

$ws.Range("H132").Value = 33338792
$ws.Range("I132").Value = 4616.1665
$ws.Range("K132").Value = 13848.4995
$ws.Range("M132").Value = -11318.4995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5623.625
$ws.Range("J3").Value = 4362
$ws.Range("L3").Value = 4362
$ws.Range("N3").Value = -4590

$ws.Range("H105").Value = 3299.5
$ws.Range("I105").Value = 3299.5
$ws.Range("K105").Value = 3299.5
$ws.Range("M105").Value = -1552.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5126.116
$ws.Range("I31").Value = 3613.2
$ws.Range("J31").Value = 5584.5757
$ws.Range("K31").Value = 3613.2
$ws.Range("L31").Value = 5584.5757
$ws.Range("M31").Value = -3318.2
$ws.Range("N31").Value = -6174.5757

$ws.Range("H34").Value = 5126.116
$ws.Range("I34").Value = 3613.2
$ws.Range("J34").Value = 5584.5757
$ws.Range("K34").Value = 3613.2
$ws.Range("L34").Value = 5584.5757
$ws.Range("M34").Value = -3411.2
$ws.Range("N34").Value = -5988.5757

$ws.Range("H105").Value = 2549.8948
$ws.Range("I105").Value = 2448.375
$ws.Range("K105").Value = 2448.375
$ws.Range("M105").Value = -701.375

$ws.Range("H122").Value = 12854.728
$ws.Range("I122").Value = 1482.6875
$ws.Range("K122").Value = 4448.0625
$ws.Range("M122").Value = -1998.0625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1481
$ws.Range("J5").Value = 1996
$ws.Range("L5").Value = 5988
$ws.Range("N5").Value = -6212

$ws.Range("H46").Value = 9093226
$ws.Range("I46").Value = 14286928
$ws.Range("J46").Value = 4247.25
$ws.Range("K46").Value = 42860784
$ws.Range("L46").Value = 12741.75
$ws.Range("M46").Value = -42860693
$ws.Range("N46").Value = -12923.75

$ws.Range("H68").Value = 3686.875
$ws.Range("I68").Value = 1708.25
$ws.Range("K68").Value = 5124.75
$ws.Range("M68").Value = -4313.75

$ws.Range("H71").Value = 3686.875
$ws.Range("I71").Value = 1708.25
$ws.Range("K71").Value = 15374.25
$ws.Range("M71").Value = -11318.25

$ws.Range("H103").Value = 447.5
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws.Range("H107").Value = 1306.5
$ws.Range("J107").Value = 1998.5714
$ws.Range("L107").Value = 5995.7142
$ws.Range("N107").Value = -9835.7142

$ws.Range("H127").Value = 1676.6
$ws.Range("I127").Value = 719.5
$ws.Range("J127").Value = 2314.6667
$ws.Range("K127").Value = 2158.5
$ws.Range("L127").Value = 6944.000100000001
$ws.Range("M127").Value = 2801.5
$ws.Range("N127").Value = -16864.0001

$ws.Range("H132").Value = 1312.1428
$ws.Range("I132").Value = 1307.8
$ws.Range("K132").Value = 11770.2
$ws.Range("M132").Value = -9240.199999999999

$ws.Range("H135").Value = 1481
$ws.Range("J135").Value = 1996
$ws.Range("L135").Value = 17964
$ws.Range("N135").Value = -23034

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3097.158
$ws.Range("I122").Value = 3646.3845
$ws.Range("J122").Value = 1907.1666
$ws.Range("K122").Value = 10939.1535
$ws.Range("L122").Value = 5721.4998
$ws.Range("M122").Value = -8489.1535
$ws.Range("N122").Value = -10621.4998

$ws.Range("H132").Value = 13582.4
$ws.Range("I132").Value = 16003
$ws.Range("K132").Value = 48009
$ws.Range("M132").Value = -45479

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4297.4
$ws.Range("J122").Value = 4498.5
$ws.Range("L122").Value = 13495.5
$ws.Range("N122").Value = -18395.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1532.8529
$ws.Range("I122").Value = 1497.4849
$ws.Range("J122").Value = 2700
$ws.Range("K122").Value = 4492.4547
$ws.Range("L122").Value = 8100
$ws.Range("M122").Value = -2042.4547
$ws.Range("N122").Value = -13000
